$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "theta_threshold_range" row (row 5) is removed entirely; rows below it
# shift up. Deleting the whole row (rather than just clearing it) re-indexes
# the remaining rows, matching the diff where what used to be row 6
# (pie_threshold_range) becomes row 5.
$ws.Rows(5).Delete()

# Update the surviving data values to their new targets.
$ws.Range("B2").Value = 5.4
$ws.Range("B3").Value = 5
$ws.Range("C4").Value = 1.4
$ws.Range("C5").Value = 15

# The old row 6's B cell (now B5) carried a special "highlight" style
# (Times New Roman, 12pt) that is no longer present on any cell after the
# edit; normalize it back to the plain data-row style used by its
# neighbours.
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the saved selection/active cell from the edited workbook.
$ws.Range("C4").Select() | Out-Null
